# Update Gpc3-Cd81 LR-pair TPM-derived NATMI metrics with recomputed values
# (ligand/receptor expression, specificity and edge-weight columns E:T)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.6654429999999999
$ws.Range("H2").Value = 1.996329
$ws.Range("I2").Value = 0.01979810471132212
$ws.Range("J2").Value = 0.01979810471132213
$ws.Range("M2").Value = 155.8020323333334
$ws.Range("N2").Value = 467.406097
$ws.Range("O2").Value = 0.376400502549802
$ws.Range("P2").Value = 0.3764005025498019
$ws.Range("Q2").Value = 103.6773718019903
$ws.Range("R2").Value = 933.096346217913
$ws.Range("S2").Value = 0.007452016562875249
$ws.Range("T2").Value = 0.007452016562875249
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.6654429999999999
$ws.Range("H3").Value = 1.996329
$ws.Range("I3").Value = 0.01979810471132212
$ws.Range("J3").Value = 0.01979810471132213
$ws.Range("O3").Value = 0.4356399183007033
$ws.Range("P3").Value = 0.4356399183007033
$ws.Range("Q3").Value = 119.994531026097
$ws.Range("R3").Value = 1079.950779234873
$ws.Range("S3").Value = 0.008624844718949138
$ws.Range("T3").Value = 0.008624844718949139
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.6654429999999999
$ws.Range("H4").Value = 1.996329
$ws.Range("I4").Value = 0.01979810471132212
$ws.Range("J4").Value = 0.01979810471132213
$ws.Range("M4").Value = 10.79259366666667
$ws.Range("N4").Value = 32.377781
$ws.Range("O4").Value = 0.02607371430982302
$ws.Range("P4").Value = 0.02607371430982302
$ws.Range("Q4").Value = 7.181855907327666
$ws.Range("R4").Value = 64.63670316594899
$ws.Range("S4").Value = 0.0005162101261189742
$ws.Range("T4").Value = 0.0005162101261189742
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 0.6654429999999999
$ws.Range("H5").Value = 1.996329
$ws.Range("I5").Value = 0.01979810471132212
$ws.Range("J5").Value = 0.01979810471132213
$ws.Range("M5").Value = 53.06048966666666
$ws.Range("N5").Value = 159.181469
$ws.Range("O5").Value = 0.1281882827647747
$ws.Range("P5").Value = 0.1281882827647747
$ws.Range("Q5").Value = 35.30873142525566
$ws.Range("R5").Value = 317.778582827301
$ws.Range("S5").Value = 0.002537885044941578
$ws.Range("T5").Value = 0.002537885044941578
# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.6654429999999999
$ws.Range("H6").Value = 1.996329
$ws.Range("I6").Value = 0.01979810471132212
$ws.Range("J6").Value = 0.01979810471132213
$ws.Range("M6").Value = 13.94831233333333
$ws.Range("N6").Value = 41.844937
$ws.Range("O6").Value = 0.03369758207489706
$ws.Range("P6").Value = 0.03369758207489706
$ws.Range("Q6").Value = 9.281806804030332
$ws.Range("R6").Value = 83.53626123627299
$ws.Range("S6").Value = 0.0006671482584371835
$ws.Range("T6").Value = 0.0006671482584371837
# Row 7
$ws.Range("I7").Value = 0.7271044512817749
$ws.Range("J7").Value = 0.727104451281775
$ws.Range("M7").Value = 155.8020323333334
$ws.Range("N7").Value = 467.406097
$ws.Range("O7").Value = 0.376400502549802
$ws.Range("P7").Value = 0.3764005025498019
$ws.Range("Q7").Value = 3807.651269331455
$ws.Range("R7").Value = 34268.86142398309
$ws.Range("S7").Value = 0.2736824808686581
$ws.Range("T7").Value = 0.2736824808686581
# Row 8
$ws.Range("I8").Value = 0.7271044512817749
$ws.Range("J8").Value = 0.727104451281775
$ws.Range("O8").Value = 0.4356399183007033
$ws.Range("P8").Value = 0.4356399183007033
$ws.Range("S8").Value = 0.3167557237524701
$ws.Range("T8").Value = 0.3167557237524701
# Row 9
$ws.Range("I9").Value = 0.7271044512817749
$ws.Range("J9").Value = 0.727104451281775
$ws.Range("M9").Value = 10.79259366666667
$ws.Range("N9").Value = 32.377781
$ws.Range("O9").Value = 0.02607371430982302
$ws.Range("P9").Value = 0.02607371430982302
$ws.Range("Q9").Value = 263.7605707629138
$ws.Range("R9").Value = 2373.845136866224
$ws.Range("S9").Value = 0.01895831373612163
$ws.Range("T9").Value = 0.01895831373612163
# Row 10
$ws.Range("I10").Value = 0.7271044512817749
$ws.Range("J10").Value = 0.727104451281775
$ws.Range("M10").Value = 53.06048966666666
$ws.Range("N10").Value = 159.181469
$ws.Range("O10").Value = 0.1281882827647747
$ws.Range("P10").Value = 0.1281882827647747
$ws.Range("Q10").Value = 1296.747146393975
$ws.Range("R10").Value = 11670.72431754578
$ws.Range("S10").Value = 0.09320627100043449
$ws.Range("T10").Value = 0.09320627100043451
# Row 11
$ws.Range("I11").Value = 0.7271044512817749
$ws.Range("J11").Value = 0.727104451281775
$ws.Range("M11").Value = 13.94831233333333
$ws.Range("N11").Value = 41.844937
$ws.Range("O11").Value = 0.03369758207489706
$ws.Range("P11").Value = 0.03369758207489706
$ws.Range("Q11").Value = 340.8832886558276
$ws.Range("R11").Value = 3067.949597902448
$ws.Range("S11").Value = 0.0245016619240906
$ws.Range("T11").Value = 0.02450166192409061
# Row 12
$ws.Range("G12").Value = 8.477506666666667
$ws.Range("H12").Value = 25.43252
$ws.Range("I12").Value = 0.2522207982916614
$ws.Range("J12").Value = 0.2522207982916615
$ws.Range("M12").Value = 155.8020323333334
$ws.Range("N12").Value = 467.406097
$ws.Range("O12").Value = 0.376400502549802
$ws.Range("P12").Value = 0.3764005025498019
$ws.Range("Q12").Value = 1320.812767786049
$ws.Range("R12").Value = 11887.31491007444
$ws.Range("S12").Value = 0.0949360352304936
$ws.Range("T12").Value = 0.0949360352304936
# Row 13
$ws.Range("G13").Value = 8.477506666666667
$ws.Range("H13").Value = 25.43252
$ws.Range("I13").Value = 0.2522207982916614
$ws.Range("J13").Value = 0.2522207982916615
$ws.Range("O13").Value = 0.4356399183007033
$ws.Range("P13").Value = 0.4356399183007033
$ws.Range("Q13").Value = 1528.68756112436
$ws.Range("R13").Value = 13758.18805011924
$ws.Range("S13").Value = 0.1098774479615175
$ws.Range("T13").Value = 0.1098774479615176
# Row 14
$ws.Range("G14").Value = 8.477506666666667
$ws.Range("H14").Value = 25.43252
$ws.Range("I14").Value = 0.2522207982916614
$ws.Range("J14").Value = 0.2522207982916615
$ws.Range("M14").Value = 10.79259366666667
$ws.Range("N14").Value = 32.377781
$ws.Range("O14").Value = 0.02607371430982302
$ws.Range("P14").Value = 0.02607371430982302
$ws.Range("Q14").Value = 91.49428475979111
$ws.Range("R14").Value = 823.44856283812
$ws.Range("S14").Value = 0.006576333037652279
$ws.Range("T14").Value = 0.00657633303765228
# Row 15
$ws.Range("G15").Value = 8.477506666666667
$ws.Range("H15").Value = 25.43252
$ws.Range("I15").Value = 0.2522207982916614
$ws.Range("J15").Value = 0.2522207982916615
$ws.Range("M15").Value = 53.06048966666666
$ws.Range("N15").Value = 159.181469
$ws.Range("O15").Value = 0.1281882827647747
$ws.Range("P15").Value = 0.1281882827647747
$ws.Range("Q15").Value = 449.8206548857644
$ws.Range("R15").Value = 4048.38589397188
$ws.Range("S15").Value = 0.03233175101056869
$ws.Range("T15").Value = 0.0323317510105687
# Row 16
$ws.Range("G16").Value = 8.477506666666667
$ws.Range("H16").Value = 25.43252
$ws.Range("I16").Value = 0.2522207982916614
$ws.Range("J16").Value = 0.2522207982916615
$ws.Range("M16").Value = 13.94831233333333
$ws.Range("N16").Value = 41.844937
$ws.Range("O16").Value = 0.03369758207489706
$ws.Range("P16").Value = 0.03369758207489706
$ws.Range("Q16").Value = 118.2469107945822
$ws.Range("R16").Value = 1064.22219715124
$ws.Range("S16").Value = 0.008499231051429319
$ws.Range("T16").Value = 0.00849923105142932
# Row 17
$ws.Range("E17").Value = 1
$ws.Range("F17").Value = 0.3333333333333333
$ws.Range("G17").Value = 0.02946533333333333
$ws.Range("H17").Value = 0.088396
$ws.Range("I17").Value = 0.0008766457152413409
$ws.Range("J17").Value = 0.000876645715241341
$ws.Range("M17").Value = 155.8020323333334
$ws.Range("N17").Value = 467.406097
$ws.Range("O17").Value = 0.376400502549802
$ws.Range("P17").Value = 0.3764005025498019
$ws.Range("Q17").Value = 4.590758816712445
$ws.Range("R17").Value = 41.31682935041201
$ws.Range("S17").Value = 0.0003299698877749713
$ws.Range("T17").Value = 0.0003299698877749713
# Row 18
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 0.3333333333333333
$ws.Range("G18").Value = 0.02946533333333333
$ws.Range("H18").Value = 0.088396
$ws.Range("I18").Value = 0.0008766457152413409
$ws.Range("J18").Value = 0.000876645715241341
$ws.Range("O18").Value = 0.4356399183007033
$ws.Range("P18").Value = 0.4356399183007033
$ws.Range("Q18").Value = 5.313270790828
$ws.Range("R18").Value = 47.819437117452
$ws.Range("S18").Value = 0.0003819018677663993
$ws.Range("T18").Value = 0.0003819018677663993
# Row 19
$ws.Range("E19").Value = 1
$ws.Range("F19").Value = 0.3333333333333333
$ws.Range("G19").Value = 0.02946533333333333
$ws.Range("H19").Value = 0.088396
$ws.Range("I19").Value = 0.0008766457152413409
$ws.Range("J19").Value = 0.000876645715241341
$ws.Range("M19").Value = 10.79259366666667
$ws.Range("N19").Value = 32.377781
$ws.Range("O19").Value = 0.02607371430982302
$ws.Range("P19").Value = 0.02607371430982302
$ws.Range("Q19").Value = 0.3180073699195556
$ws.Range("R19").Value = 2.862066329276
$ws.Range("S19").Value = 0.00002285740993013319
$ws.Range("T19").Value = 0.00002285740993013319
# Row 20
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 0.3333333333333333
$ws.Range("G20").Value = 0.02946533333333333
$ws.Range("H20").Value = 0.088396
$ws.Range("I20").Value = 0.0008766457152413409
$ws.Range("J20").Value = 0.000876645715241341
$ws.Range("M20").Value = 53.06048966666666
$ws.Range("N20").Value = 159.181469
$ws.Range("O20").Value = 0.1281882827647747
$ws.Range("P20").Value = 0.1281882827647747
$ws.Range("Q20").Value = 1.563445014858222
$ws.Range("R20").Value = 14.071005133724
$ws.Range("S20").Value = 0.0001123757088298851
$ws.Range("T20").Value = 0.0001123757088298852
# Row 21
$ws.Range("E21").Value = 1
$ws.Range("F21").Value = 0.3333333333333333
$ws.Range("G21").Value = 0.02946533333333333
$ws.Range("H21").Value = 0.088396
$ws.Range("I21").Value = 0.0008766457152413409
$ws.Range("J21").Value = 0.000876645715241341
$ws.Range("M21").Value = 13.94831233333333
$ws.Range("N21").Value = 41.844937
$ws.Range("O21").Value = 0.03369758207489706
$ws.Range("P21").Value = 0.03369758207489706
$ws.Range("Q21").Value = 0.4109916723391111
$ws.Range("R21").Value = 3.698925051052
$ws.Range("S21").Value = 0.00002954084093995193
$ws.Range("T21").Value = 0.00002954084093995193
